# Update "want to go" (想去人数) counts for several events across sheets,
# matching the data refresh described in the commit message
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 12984
$ws1.Range("F5").Value  = 81
$ws1.Range("F6").Value  = 92
$ws1.Range("F7").Value  = 51
$ws1.Range("F8").Value  = 23
$ws1.Range("F10").Value = 12964
$ws1.Range("F11").Value = 292
$ws1.Range("F12").Value = 44
$ws1.Range("F13").Value = 8714
$ws1.Range("F14").Value = 7723
$ws1.Range("F15").Value = 204
$ws1.Range("F16").Value = 116
$ws1.Range("F18").Value = 130
$ws1.Range("F20").Value = 11
$ws1.Range("F24").Value = 328
$ws1.Range("F26").Value = 5219

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 12984
$ws4.Range("F6").Value  = 81
$ws4.Range("F7").Value  = 92
$ws4.Range("F8").Value  = 51
$ws4.Range("F9").Value  = 23
$ws4.Range("F11").Value = 12964
$ws4.Range("F12").Value = 292
$ws4.Range("F13").Value = 44
$ws4.Range("F14").Value = 8714
$ws4.Range("F15").Value = 7723
$ws4.Range("F16").Value = 204
$ws4.Range("F17").Value = 116
$ws4.Range("F19").Value = 130
$ws4.Range("F21").Value = 11
$ws4.Range("F23").Value = 1
$ws4.Range("F27").Value = 328
$ws4.Range("F29").Value = 5219
